# Latest update for the documentation.
#
# The paragraph ending in "... Model as a string ... end user." currently
# also carries the (invisible) "_GoBack" bookmark right at its end,
# immediately followed by one empty trailing paragraph.
#
# Target shape:
#   ... end user.                                  (unchanged text, bookmark removed from it)
#   Ordering and sorting are the same thing         (new paragraph)
#   <empty paragraph containing just the bookmark>  (bookmark now lives alone)
#   <empty paragraph>                               (new)
#   <empty paragraph>                               (pre-existing trailing paragraph)

$d = $word.ActiveDocument

# Step 1: split the existing paragraph right after "... end user." into
# three paragraphs: the original sentence, the new "Ordering and
# sorting..." sentence, and a paragraph break before whatever trails the
# matched text (the _GoBack bookmark). Word keeps content that trailed the
# match attached after the final inserted paragraph mark, so this leaves
# the bookmark alone in its own new paragraph - exactly like the target.
$find = $d.Content.Find
$oldText = "more convenient for the end user."
$newText = "more convenient for the end user.^pOrdering and sorting are the same thing^p"
$find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# Step 2: locate the paragraph we just created ("Ordering and sorting...")
# and then insert one more empty paragraph right after the paragraph that
# follows it (the one now holding only the bookmark), before the
# already-existing trailing empty paragraph.
$orderingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($paraText -eq "Ordering and sorting are the same thing") {
        $orderingIndex = $i
        break
    }
}

$bookmarkPara = $d.Paragraphs.Item($orderingIndex + 1)
$r = $bookmarkPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
